# Update countries & provincias Spain
# Applies the refreshed COVID data snapshot (Datos actualizados ... 03:22)
# to the "Pais" worksheet: updates the footer timestamp, refreshes the
# numeric counters for several countries, and swaps two pairs of country
# names (Reunion/Jamaica and Guadalupe/Birmania) whose rows kept their
# position but exchanged labels + data as the new snapshot re-ranked them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (row 1, A1) ---------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 03:22"

# --- Row 4: Estados Unidos ------------------------------------------------
$ws.Range("B4").Value = 1095023
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 152324
$ws.Range("E4").Value = 878843
$ws.Range("F4").Value = 15226
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 63856

# --- Row 13: Brasil --------------------------------------------------------
$ws.Range("B13").Value = 87187
$ws.Range("C13").Value = 1807
$ws.Range("D13").Value = 35935
$ws.Range("E13").Value = 45246
$ws.Range("F13").Value = 8318
$ws.Range("G13").Value = 105
$ws.Range("H13").Value = 6006

# --- Row 48: Australia -------------------------------------------------
$ws.Range("B48").Value = 6766
$ws.Range("C48").Value = 12
$ws.Range("D48").Value = 5739
$ws.Range("E48").Value = 934
$ws.Range("F48").Value = 31
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 93

# --- Row 82: Nueva Zelanda ------------------------------------------------
$ws.Range("B82").Value = 1479
$ws.Range("C82").Value = 3
$ws.Range("D82").Value = 1252
$ws.Range("E82").Value = 208
$ws.Range("F82").Value = 1
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 19

# --- Rows 119/120: Reunion and Jamaica swap labels + refreshed data -------
$ws.Range("A119").Value = "Jamaica"
$ws.Range("B119").Value = 422
$ws.Range("C119").Value = 26
$ws.Range("D119").Value = 29
$ws.Range("E119").Value = 385
$ws.Range("F119").Value = 3
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 8

$ws.Range("A120").Value = "Reunion"
$ws.Range("B120").Value = 420
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 300
$ws.Range("E120").Value = 120
$ws.Range("F120").Value = 2
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 0

# --- Rows 137/138: Guadalupe and Birmania swap labels + refreshed data ----
$ws.Range("A137").Value = "Birmania"
$ws.Range("B137").Value = 151
$ws.Range("C137").Value = 1
$ws.Range("D137").Value = 27
$ws.Range("E137").Value = 118
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 6

$ws.Range("A138").Value = "Guadalupe"
$ws.Range("B138").Value = 151
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 95
$ws.Range("E138").Value = 44
$ws.Range("F138").Value = 11
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 12
